# "Aanvullen van de tijdsbesteding"
# Add a "gemiddeld:" (average) row below the existing "totaal:" (total) row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New label in D3, matching the "totaal:" label already in D2.
$ws.Range("D3").Value = "gemiddeld:"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats - copy formatting only

# New average formula in E3, matching the formatting used by the SUM in E2.
$ws.Range("E3").Formula = "=AVERAGE(B2:B8)"
$ws.Range("E2").Copy() | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats - copy formatting only

$excel.CutCopyMode = $false

# Move the active selection to B4.
$ws.Range("B4").Select() | Out-Null
